$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 54-83: "Electives"
for ($r = 54; $r -le 83; $r++) {
    $ws.Cells.Item($r, 4).Value = "Electives"
}

# Rows 84-89: "Other Requirements"
for ($r = 84; $r -le 89; $r++) {
    $ws.Cells.Item($r, 4).Value = "Other Requirements"
}
